# Revert the "Artificial Intelligence and Machine Learning" themed template
# text back to the original "Information Technology" budget template wording.
# (Reverts commit 34e0294 "RESTORE: Recover all 973 original multi-industry
# template files from commit 168d9c4" for this workbook.)

$wb = $excel.ActiveWorkbook

# --- Instructions & User Guide -------------------------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")
$ws.Range("A1").Value = "Information Technology Comprehensive Budget - User Guide & Instructions"
$ws.Range("A56").Value = "📋 INFORMATION TECHNOLOGY PROJECT OVERVIEW"
$ws.Range("B59").Value = "IT Managers, DevOps Engineers, AI Architects, DevOps Engineers..."

# --- Budget Summary --------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Range("A1").Value = "Information Technology - Executive Budget Summary"

# --- Resources ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")
$ws.Range("A1").Value = "Information Technology - Resources Budget"
$ws.Range("A4").Value = "IT Managers"
$ws.Range("A5").Value = "DevOps Engineers"
$ws.Range("A9").Value = "System Administrators"

# --- Logistics -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")
$ws.Range("A1").Value = "Information Technology - Logistics Budget"

# --- Technology ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")
$ws.Range("A1").Value = "Information Technology - Technology Budget"

# --- Training ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")
$ws.Range("A1").Value = "Information Technology - Training Budget"
$ws.Range("A4").Value = "IT Certification Programs"

# --- Contingency --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")
$ws.Range("A1").Value = "Information Technology - Contingency Budget"

# --- Timeline -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")
$ws.Range("A1").Value = "Information Technology - Budget Timeline"
